$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts old E:T -> F:U)
$ws.Columns("E:E").Insert()

# Header for the new column
$ws.Range("E1").Value = "UBID"

# New column has no special formatting (reset any inherited style)
$ws.Range("E2:E15").Style = "Normal"

# Fill in the UBID values for each data row
$ws.Range("E2").Value = "https://ubid.org/65923-510"
$ws.Range("E3").Value = "https://ubid.org/64942-1191"
$ws.Range("E4").Value = "https://ubid.org/0093-2068"
$ws.Range("E5").Value = "https://ubid.org/55316-267"
$ws.Range("E6").Value = "https://ubid.org/62011-0243"
$ws.Range("E7").Value = "https://ubid.org/49999-504"
$ws.Range("E8").Value = "https://ubid.org/55289-460"
$ws.Range("E9").Value = "https://ubid.org/60505-3404"
$ws.Range("E10").Value = "https://ubid.org/55319-377"
$ws.Range("E11").Value = "https://ubid.org/55505-167"
$ws.Range("E12").Value = "https://ubid.org/61062-0007"
$ws.Range("E13").Value = "https://ubid.org/0268-0851"
$ws.Range("E14").Value = "https://ubid.org/68151-1305"
$ws.Range("E15").Value = "https://ubid.org/54575-933"

# Match the final selection state
$ws.Range("E10").Select()
